# Commit: "Append: 2025-12-16 18:37 JST"
#
# The scraper re-ran and refreshed the acquisition timestamp stored in
# column A ("取得日時") of the active "ランサーズ" sheet for every existing
# data row (rows 2-15), bumping it from 2025-12-16 18:28:50 to
# 2025-12-16 18:37:40. No other cell values change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-12-16 18:37:40"

$firstDataRow = 2
$lastDataRow = 15

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
